# Update countries & provincias Spain
# - Reorders three country pairs in the ranking (because their case counts
#   changed position relative to their neighbour)
# - Refreshes the "Datos actualizados" timestamp
# - Refreshes the numeric COVID-19 stats for the affected rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name swaps (rows stay sorted by "Casos totales") ---------
# Filipinas / Colombia / Dinamarca / Serbia block: Colombia now ranks
# above Dinamarca.
$ws.Range("A43").Value = "Colombia"
$ws.Range("A44").Value = "Dinamarca"

# Ghana / Nigeria / Afganistan / Luxemburgo block: Nigeria now ranks
# above Afganistan.
$ws.Range("A63").Value = "Nigeria"
$ws.Range("A64").Value = "Afganistan"

# Papua Nueva Guinea / Butan / Islas Virgenes Britanicas block: Butan now
# ranks above Islas Virgenes Britanicas.
$ws.Range("A212").Value = "Butan"
$ws.Range("A213").Value = "Islas Virgenes Britanicas"

# --- Timestamp update ---------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 10 de Mayo de 2020 a las 01:04"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos, -
#     Recuperados, Casos criticos, Muertes hoy, Muertes) -----------------

# Estados Unidos
$ws.Range("B4").Value = 1346339
$ws.Range("C4").Value = 24554
$ws.Range("D4").Value = 236825
$ws.Range("E4").Value = 1029510
$ws.Range("F4").Value = 16811
$ws.Range("G4").Value = 1389
$ws.Range("H4").Value = 80004

# Brasil
$ws.Range("B11").Value = 155939
$ws.Range("C11").Value = 10047
$ws.Range("E11").Value = 86015
$ws.Range("G11").Value = 635
$ws.Range("H11").Value = 10627

# Canada
$ws.Range("B15").Value = 67702
$ws.Range("C15").Value = 1268
$ws.Range("D15").Value = 31249
$ws.Range("E15").Value = 31760
$ws.Range("G15").Value = 124
$ws.Range("H15").Value = 4693

# Suiza
$ws.Range("D22").Value = 26400
$ws.Range("E22").Value = 2021

# Colombia (now row 43)
$ws.Range("B43").Value = 10495
$ws.Range("C43").Value = 444
$ws.Range("D43").Value = 2569
$ws.Range("E43").Value = 7481
$ws.Range("F43").Value = 130
$ws.Range("G43").Value = 17
$ws.Range("H43").Value = 445

# Dinamarca (now row 44)
$ws.Range("B44").Value = 10319
$ws.Range("C44").Value = 101
$ws.Range("D44").Value = 8093
$ws.Range("E44").Value = 1700
$ws.Range("F44").Value = 39
$ws.Range("G44").Value = 4
$ws.Range("H44").Value = 526

# Argentina
$ws.Range("B57").Value = 5766
$ws.Range("C57").Value = 155
$ws.Range("E57").Value = 3738
$ws.Range("G57").Value = 7
$ws.Range("H57").Value = 300

# Nigeria (now row 63)
$ws.Range("B63").Value = 4151
$ws.Range("C63").Value = 239
$ws.Range("D63").Value = 745
$ws.Range("E63").Value = 3278
$ws.Range("F63").Value = 4
$ws.Range("G63").Value = 11
$ws.Range("H63").Value = 128

# Afganistan (now row 64)
$ws.Range("B64").Value = 4033
$ws.Range("C64").Value = 255
$ws.Range("D64").Value = 502
$ws.Range("E64").Value = 3416
$ws.Range("F64").Value = 7
$ws.Range("G64").Value = 6
$ws.Range("H64").Value = 115

# Sudan
$ws.Range("B92").Value = 1164
$ws.Range("C92").Value = 53
$ws.Range("D92").Value = 119
$ws.Range("E92").Value = 981
$ws.Range("G92").Value = 5
$ws.Range("H92").Value = 64

# Maldivas
$ws.Range("B106").Value = 790
$ws.Range("C106").Value = 46
$ws.Range("E106").Value = 758

# Uruguay
$ws.Range("B112").Value = 702
$ws.Range("C112").Value = 8
$ws.Range("D112").Value = 513
$ws.Range("E112").Value = 171
$ws.Range("F112").Value = 8

# Monaco
$ws.Range("B163").Value = 96
$ws.Range("C163").Value = 1
$ws.Range("E163").Value = 10

# Butan (now row 212)
$ws.Range("D212").Value = 5
$ws.Range("H212").Value = 0

# Islas Virgenes Britanicas (now row 213)
$ws.Range("D213").Value = 4
$ws.Range("H213").Value = 1
